$d = $word.ActiveDocument

# Paragraph 1: "Cd  Change Dic" -> "Ls List dic"
$p1 = $d.Paragraphs.Item(1).Range
$p1.Find.Execute("Cd  Change Dic", $false, $false, $false, $false, $false, `
                  $true, 1, $false, "Ls List dic", 2)

# Paragraph 2: "Ls List dic" -> "Pwd present working dic"
$p2 = $d.Paragraphs.Item(2).Range
$p2.Find.Execute("Ls List dic", $false, $false, $false, $false, $false, `
                  $true, 1, $false, "Pwd present working dic", 2)

# Paragraph 3: "Pwd present working dic" -> "Mdkir "
$p3 = $d.Paragraphs.Item(3).Range
$p3.Find.Execute("Pwd present working dic", $false, $false, $false, $false, $false, `
                  $true, 1, $false, "Mdkir ", 2)

# Paragraph 4 ("Mdkir ") is now redundant -- remove it (and its paragraph mark)
# entirely so the following empty paragraphs shift up, matching the diff.
$p4 = $d.Paragraphs.Item(4).Range
$p4.Delete()

# Last paragraph: "Add somthing" -> "Do  somthing"
$last = $d.Paragraphs.Last.Range
$last.Find.Execute("Add somthing", $false, $false, $false, $false, $false, `
                    $true, 1, $false, "Do  somthing", 2)
